$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the price/volume columns so that numeric-looking
# strings (e.g. "44.80", "1.000", "28.315.26") are stored as literal text and
# are not auto-coerced into floating point numbers (which would lose formatting
# such as trailing zeros or thousand-separator dots).
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range('D2').Value2 = '28.315.26'
$ws.Range('E2').Value2 = '  +2.54%  '
$ws.Range('D3').Value2 = '1.813.43'
$ws.Range('E3').Value2 = '  +3.74%  '
$ws.Range('D4').Value2 = '1.001'
$ws.Range('E4').Value2 = '  +0.15%  '
$ws.Range('D5').Value2 = '326.96'
$ws.Range('E5').Value2 = '  +1.35%  '
$ws.Range('D6').Value2 = '0.9999'
$ws.Range('E6').Value2 = '  +0.16%  '
$ws.Range('D7').Value2 = '0.4365'
$ws.Range('E7').Value2 = '  +2.50%  '
$ws.Range('D8').Value2 = '0.3672'
$ws.Range('E8').Value2 = '  +1.52%  '
$ws.Range('D9').Value2 = '44.80'
$ws.Range('E9').Value2 = '  -1.19%  '
$ws.Range('D10').Value2 = '0.07685'
$ws.Range('E10').Value2 = '  +2.79%  '
$ws.Range('D11').Value2 = '1.143'
$ws.Range('D12').Value2 = '1.000'
$ws.Range('E12').Value2 = '  +0.21%  '
$ws.Range('D13').Value2 = '22.05'
$ws.Range('E13').Value2 = '  +2.39%  '
$ws.Range('D14').Value2 = '6.313'
$ws.Range('E14').Value2 = '  +3.30%  '
$ws.Range('D15').Value2 = '7.529'
$ws.Range('E15').Value2 = '  +4.26%  '
$ws.Range('D16').Value2 = '1.823.79'
$ws.Range('E16').Value2 = '  +4.73%  '
$ws.Range('D17').Value2 = '95.68'
$ws.Range('E17').Value2 = '  +9.13%  '
$ws.Range('D18').Value2 = '0.00001081'
$ws.Range('E18').Value2 = '  +1.16%  '
$ws.Range('D19').Value2 = '0.06526'
$ws.Range('E19').Value2 = '  +4.63%  '
$ws.Range('D20').Value2 = '0.9992'
$ws.Range('E20').Value2 = '  +0.12%  '
$ws.Range('E21').Value2 = '  +2.83%  '
$ws.Range('D22').Value2 = '6.243'
$ws.Range('E22').Value2 = '  +1.85%  '
$ws.Range('D23').Value2 = '28.323.05'
$ws.Range('E23').Value2 = '  +2.63%  '
$ws.Range('D24').Value2 = '11.59'
$ws.Range('E24').Value2 = '  -0.55%  '
$ws.Range('D25').Value2 = '2.084'
$ws.Range('E25').Value2 = '  -10.02%  '
$ws.Range('D26').Value2 = '162.30'
$ws.Range('E26').Value2 = '  +7.02%  '
$ws.Range('E27').Value2 = '  +1.19%  '
$ws.Range('D28').Value2 = '2.026.68'
$ws.Range('E28').Value2 = '  +4.49%  '
$ws.Range('E29').Value2 = '  -3.91%  '
$ws.Range('D30').Value2 = '128.99'
$ws.Range('E30').Value2 = '  +1.80%  '
$ws.Range('D31').Value2 = '1.208'
$ws.Range('E31').Value2 = '  -0.99%  '
$ws.Range('D32').Value2 = '5.963'
$ws.Range('E32').Value2 = '  +4.48%  '
$ws.Range('D33').Value2 = '0.09192'
$ws.Range('E33').Value2 = '  +0.44%  '
$ws.Range('D34').Value2 = '3.500'
$ws.Range('E34').Value2 = '  -4.86%  '
$ws.Range('E35').Value2 = '  +2.41%  '
$ws.Range('D36').Value2 = '0.02346'
$ws.Range('E36').Value2 = '  +1.78%  '
$ws.Range('D37').Value2 = '5.198'
$ws.Range('E37').Value2 = '  +2.15%  '
$ws.Range('D38').Value2 = '0.2171'
$ws.Range('E38').Value2 = '  +1.57%  '
$ws.Range('D39').Value2 = '0.6592'
$ws.Range('E39').Value2 = '  +2.35%  '
$ws.Range('D40').Value2 = '0.06213'
$ws.Range('E40').Value2 = '  +1.86%  '
$ws.Range('D41').Value2 = '1.196'
$ws.Range('E41').Value2 = '  +0.10%  '
$ws.Range('D42').Value2 = '8.132'
$ws.Range('E42').Value2 = '  +2.46%  '
$ws.Range('E43').Value2 = '  +1.11%  '
$ws.Range('D44').Value2 = '0.9994'
$ws.Range('E44').Value2 = '  +0.17%  '
$ws.Range('D45').Value2 = '13.97'
$ws.Range('E45').Value2 = '  +2.04%  '
$ws.Range('D46').Value2 = '0.6116'
$ws.Range('E46').Value2 = '  +3.54%  '
$ws.Range('D47').Value2 = '3.747'
$ws.Range('E47').Value2 = '  +0.57%  '
$ws.Range('D48').Value2 = '125.90'
$ws.Range('E48').Value2 = '  +0.03%  '
$ws.Range('E49').Value2 = '  +2.75%  '
$ws.Range('E50').Value2 = '  +3.11%  '
$ws.Range('D51').Value2 = '0.07002'
$ws.Range('E51').Value2 = '  +1.91%  '

# Restore the default (General) style on these cells so the saved workbook
# does not leave a lingering explicit style index on cells that originally had
# none, matching the source formatting exactly.
$priceVolumeRange.Style = "Normal"

